$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A slightly (13.85546875 -> 14.85546875)
$ws.Columns.Item(1).ColumnWidth = 14

$data = @(
    @(42600.792256944442, "Random", 0, 0, 0, 0, 0, 53, 47, 0, 0, 64, 36),
    @(42600.794583333336, "Random", 0, 0, 0, 0, 0, 70, 30, 0, 0, 76, 24),
    @(42600.830914351849, "Random", 0, 0, 0, 0, 0, 8, 92, 0, 0, 82, 18),
    @(42600.879363425927, "Random", 0, 0, 0, 0, 0, 27, 73, 0, 0, 22, 78)
)

$row = 3
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    for ($col = 3; $col -le 13; $col++) {
        $ws.Cells.Item($row, $col).Value = $rec[$col - 1]
    }
    $row++
}

# Reuse the same date/time style already applied to A2 for the new A column cells
$ws.Range("A2").Copy()
$ws.Range("A3:A6").PasteSpecial(-4122) # xlPasteFormats

